$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column D (old "Terms Typically Offered"
# column), shifting it to column G, to make room for the new
# Corequisites / Concurrent / Recommended columns.
$ws.Range("D1:F1").EntireColumn.Insert()

# New header row values for the inserted columns.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Fill the new data cells (rows 2-5) with "NA".
$ws.Range("D2:F5").Value = "NA"
